$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# The "Micheal" entries in column A (rows 2-8) are actually "wiktor" -
# fix the data so column A reads "wiktor" for every data row.
$ws.Range("A2:A8").Value = "wiktor"

# Reflect the new selection left behind in the saved file: A2:A8 selected,
# with A2 as the active cell.
$ws.Activate()
$null = $ws.Range("A2:A8").Select()
